$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of cell address -> new value, per the recorded diff (case with 380 kV).
$cellValues = @{
    "B2" = 1.328088799248235; "C2" = 0.3005915337488432; "D2" = 0.02691717614016298; "F2" = 1.26782495163323; "G2" = 0.002432988663422541
    "B3" = 1.193643105224794; "C3" = 0.262740149758713; "D3" = 0.02753636415298999; "F3" = 1.228638135340702; "G3" = 0.002438568569355478
    "B4" = 1.111593090356678; "C4" = 0.2395215065747891; "D4" = 0.0279428179114749; "F4" = 1.205743070340091; "G4" = 0.002442170173400411
    "B5" = 1.07828126715691; "C5" = 0.2300648786759609; "D5" = 0.02811497430915111; "F5" = 1.196702219728536; "G5" = 0.002443682151862876
    "B6" = 1.072757323128542; "C6" = 0.2284949118759698; "D6" = 0.02814395232624989; "F6" = 1.195218333684878; "G6" = 0.002443935894596592
    "B7" = 1.111143334485291; "C7" = 0.2393939508685037; "D7" = 0.02794511336970196; "F7" = 1.205619977192271; "G7" = 0.002442190384804346
    "B8" = 1.281627116141408; "C8" = 0.2875354753128079; "D8" = 0.02712517013160287; "F8" = 1.254068908013267; "G8" = 0.002434876282431737
    "B9" = 1.620004439777802; "C9" = 0.3821427970212312; "D9" = 0.02572939304608468; "F9" = 1.358516464262038; "G9" = 0.002421918690452449
    "B10" = 1.871239517756464; "C10" = 0.4518168925330315; "D10" = 0.02483829542030058; "F10" = 1.441278503134072; "G10" = 0.002413232988602854
    "B11" = 1.986140005250547; "C11" = 0.4835595720141441; "D11" = 0.02446319247398421; "F11" = 1.48029611276641; "G11" = 0.002409460578184491
    "B12" = 2.029740373998266; "C12" = 0.4955872380373307; "D12" = 0.02432559549528968; "F12" = 1.495272534411953; "G12" = 0.002408057600583539
    "B13" = 2.020346217090719; "C13" = 0.4929965265847045; "D13" = 0.02435503019190222; "F13" = 1.492038069591729; "G13" = 0.002408358622772135
    "B14" = 1.989725217470948; "C14" = 0.4845489436963817; "D14" = 0.02445178262867032; "F14" = 1.481524171897547; "G14" = 0.002409344643124601
    "B15" = 1.970980758251812; "C15" = 0.4793755362530874; "D15" = 0.02451162817021135; "F15" = 1.475110454493091; "G15" = 0.002409951932201507
    "B16" = 1.863742991564095; "C16" = 0.4497434386053669; "D16" = 0.02486342629915761; "F16" = 1.43875653514678; "G16" = 0.002413483108615777
    "B17" = 1.798114331581246; "C17" = 0.4315776848216046; "D17" = 0.02508706090430124; "F17" = 1.416808326593667; "G17" = 0.002415695046062581
    "B18" = 1.760423979538189; "C18" = 0.4211336676367523; "D18" = 0.0252185322696743; "F18" = 1.404312860797418; "G18" = 0.002416984129146196
    "B19" = 1.747672493373045; "C19" = 0.4175982494487016; "D19" = 0.02526353157433014; "F19" = 1.40010405572626; "G19" = 0.002417423486480366
    "B20" = 1.805094647504461; "C20" = 0.4335109964128492; "D20" = 0.02506295973206818; "F20" = 1.419131412728035; "G20" = 0.002415457840283029
    "B21" = 1.998716887221917; "C21" = 0.4870299980938171; "D21" = 0.02442324265864038; "F21" = 1.484606859473274; "G21" = 0.002409054332810512
    "B22" = 2.125786060306609; "C22" = 0.5220511815921896; "D22" = 0.02403111216542442; "F22" = 1.528574276568378; "G22" = 0.002405018138138759
    "B23" = 2.057918058970188; "C23" = 0.5033555527725753; "D23" = 0.02423799250254177; "F23" = 1.504999000003266; "G23" = 0.002407158760324947
    "B24" = 1.801938720866474; "C24" = 0.4326369471673388; "D24" = 0.02507384684404101; "F24" = 1.418080763169314; "G24" = 0.002415565027075453
    "B25" = 1.528014137225512; "C25" = 0.3565230283650749; "D25" = 0.02608375848551781; "F25" = 1.329222269488668; "G25" = 0.002425276818243471
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
